$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "276.32"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.15"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.257"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06241"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.552"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.534"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "6.556"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8253"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1653"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08275"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03488"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03165"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09146"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.756"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001644"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04680"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006237"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006221"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001067"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001497"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.01397"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3292"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1249"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002731"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04750"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005290"
$ws.Range("E41").Value = "40CEJICEJI"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007017"
$ws.Range("E42").Value = "41KickTokenKICKBestin24h"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1119"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01131"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006335"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7216"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001394"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00001896"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01238"
